# DailyRoutine.xlsx -- "Questions solved before 31 July."
#
# Fills in the daily LeetCode tracker rows logged between 22 Jun 2023
# (row 12) and 20 Jul 2023 (row 21), and adds the new Easy/Medium/Hard
# daily-count columns (U/V/W) together with their header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: new Easy / Medium / Hard daily-count columns -------------
$ws.Range("U1").Value = "E"
$ws.Range("V1").Value = "M"
$ws.Range("W1").Value = "H"

# Format-only copy/paste from an existing date cell (A2) so every new date
# cell in column A reuses the workbook's existing date cell style instead
# of minting a new one.
$ws.Range("A2").Copy() | Out-Null

# Row 12 -- 22 Jun 2023
$ws.Range("A12").Value = 45099
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").Value = "+"
$ws.Range("C12").Value = "+"
$ws.Range("D12").Value = "+"
$ws.Range("E12").Value = "+"
$ws.Range("F12").Value = "+"
$ws.Range("G12").Value = "+"
$ws.Range("H12").Value = "+"
$ws.Range("I12").Value = "+"
$ws.Range("J12").Value = "+"
$ws.Range("K12").Value = "-"
$ws.Range("L12").Value = "+"
$ws.Range("M12").Value = "+"
$ws.Range("N12").Value = "+"
$ws.Range("O12").Value = "+"
$ws.Range("P12").Value = "+"
$ws.Range("Q12").Value = "+"
$ws.Range("R12").Value = "+"
$ws.Range("S12").Value = "+"
$ws.Range("T12").Value = "+"

# Row 14 -- 1 Jul 2023
$ws.Range("A14").Value = 45108
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = "N"
$ws.Range("C14").Value = "N"
$ws.Range("D14").Value = "Y"
$ws.Range("E14").Value = "N"
$ws.Range("F14").Value = "N"
$ws.Range("G14").Value = "N"
$ws.Range("H14").Value = "Y"
$ws.Range("I14").Value = "Y"
$ws.Range("J14").Value = "N"
$ws.Range("K14").Value = "N"
$ws.Range("L14").Value = "Y"
$ws.Range("M14").Value = "N"
$ws.Range("N14").Value = "N"
$ws.Range("O14").Value = "Y"
$ws.Range("P14").Value = "N"
$ws.Range("Q14").Value = "Y"
$ws.Range("R14").Value = "Y"
$ws.Range("S14").Value = "Y"
$ws.Range("T14").Value = "Y"
$ws.Range("U14").Value = 0
$ws.Range("V14").Value = 1
$ws.Range("W14").Value = 0

# Row 15 -- 14 Jul 2023
$ws.Range("A15").Value = 45121
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = "Y"
$ws.Range("C15").Value = "Y"
$ws.Range("D15").Value = "Y"
$ws.Range("E15").Value = "Y"
$ws.Range("F15").Value = "N"
$ws.Range("G15").Value = "N"
$ws.Range("H15").Value = "N"
$ws.Range("I15").Value = "N"
$ws.Range("J15").Value = "N"
$ws.Range("K15").Value = "N"
$ws.Range("L15").Value = "N"
$ws.Range("M15").Value = "Y"
$ws.Range("N15").Value = "Y"
$ws.Range("O15").Value = "Y"
$ws.Range("P15").Value = "Y"
$ws.Range("Q15").Value = "N"
$ws.Range("R15").Value = "Y"
$ws.Range("S15").Value = "Y"
$ws.Range("T15").Value = "Y"
$ws.Range("U15").Value = 3
$ws.Range("V15").Value = 0
$ws.Range("W15").Value = 0

# Row 16 -- 15 Jul 2023
$ws.Range("A16").Value = 45122
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("U16").Value = 0
$ws.Range("V16").Value = 1
$ws.Range("W16").Value = 0

# Row 17 -- 16 Jul 2023
$ws.Range("A17").Value = 45123
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = "Y"
$ws.Range("C17").Value = "Y"
$ws.Range("D17").Value = "N"
$ws.Range("E17").Value = "N"
$ws.Range("F17").Value = "N"
$ws.Range("G17").Value = "N"
$ws.Range("H17").Value = "N"
$ws.Range("I17").Value = "N"
$ws.Range("J17").Value = "N"
$ws.Range("K17").Value = "N"
$ws.Range("L17").Value = "N"
$ws.Range("M17").Value = "N"
$ws.Range("N17").Value = "N"
$ws.Range("O17").Value = "Y"
$ws.Range("P17").Value = "N"
$ws.Range("Q17").Value = "N"
$ws.Range("R17").Value = "Y"
$ws.Range("S17").Value = "Y"
$ws.Range("T17").Value = "Y"
$ws.Range("U17").Value = 1
$ws.Range("V17").Value = 1
$ws.Range("W17").Value = 0

# Row 18 -- 17 Jul 2023
$ws.Range("A18").Value = 45124
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = "Y"
$ws.Range("C18").Value = "Y"
$ws.Range("D18").Value = "Y"
$ws.Range("E18").Value = "N"
$ws.Range("F18").Value = "N"
$ws.Range("G18").Value = "N"
$ws.Range("H18").Value = "N"
$ws.Range("I18").Value = "N"
$ws.Range("J18").Value = "N"
$ws.Range("K18").Value = "N"
$ws.Range("L18").Value = "N"
$ws.Range("M18").Value = "N"
$ws.Range("N18").Value = "N"
$ws.Range("O18").Value = "N"
$ws.Range("P18").Value = "N"
$ws.Range("Q18").Value = "N"
$ws.Range("R18").Value = "Y"
$ws.Range("S18").Value = "Y"
$ws.Range("T18").Value = "Y"
$ws.Range("U18").Value = 0
$ws.Range("V18").Value = 0
$ws.Range("W18").Value = 0

# Row 19 -- 18 Jul 2023
$ws.Range("A19").Value = 45125
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").Value = "N"
$ws.Range("C19").Value = "N"
$ws.Range("D19").Value = "N"
$ws.Range("E19").Value = "Y"
$ws.Range("F19").Value = "Y"
$ws.Range("G19").Value = "N"
$ws.Range("H19").Value = "N"
$ws.Range("I19").Value = "N"
$ws.Range("J19").Value = "N"
$ws.Range("K19").Value = "N"
$ws.Range("L19").Value = "Y"
$ws.Range("M19").Value = "Y"
$ws.Range("N19").Value = "N"
$ws.Range("O19").Value = "N"
$ws.Range("P19").Value = "N"
$ws.Range("Q19").Value = "N"
$ws.Range("R19").Value = "N"
$ws.Range("S19").Value = "Y"
$ws.Range("T19").Value = "Y"
$ws.Range("U19").Value = 0
$ws.Range("V19").Value = 1
$ws.Range("W19").Value = 0

# Row 20 -- 19 Jul 2023
$ws.Range("A20").Value = 45126
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("B20").Value = "N"
$ws.Range("C20").Value = "Y"
$ws.Range("D20").Value = "Y"
$ws.Range("E20").Value = "Y"
$ws.Range("F20").Value = "N"
$ws.Range("G20").Value = "N"
$ws.Range("H20").Value = "N"
$ws.Range("I20").Value = "N"
$ws.Range("J20").Value = "N"
$ws.Range("K20").Value = "N"
$ws.Range("L20").Value = "N"
$ws.Range("M20").Value = "Y"
$ws.Range("N20").Value = "N"
$ws.Range("O20").Value = "N"
$ws.Range("P20").Value = "N"
$ws.Range("Q20").Value = "N"
$ws.Range("R20").Value = "N"
$ws.Range("S20").Value = "Y"
$ws.Range("T20").Value = "Y"
$ws.Range("U20").Value = 0
$ws.Range("V20").Value = 1
$ws.Range("W20").Value = 0

# Row 21 -- 20 Jul 2023
$ws.Range("A21").Value = 45127
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Value = "N"
$ws.Range("C21").Value = "N"
$ws.Range("D21").Value = "N"
$ws.Range("E21").Value = "N"
$ws.Range("F21").Value = "N"
$ws.Range("G21").Value = "N"
$ws.Range("H21").Value = "N"
$ws.Range("I21").Value = "N"
$ws.Range("J21").Value = "N"
$ws.Range("K21").Value = "N"
$ws.Range("L21").Value = "N"
$ws.Range("M21").Value = "Y"
$ws.Range("N21").Value = "N"
$ws.Range("O21").Value = "N"
$ws.Range("P21").Value = "N"
$ws.Range("Q21").Value = "N"
$ws.Range("R21").Value = "Y"
$ws.Range("S21").Value = "Y"
$ws.Range("T21").Value = "Y"
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = 1
$ws.Range("W21").Value = 0

# Restore the workbook's on-screen selection to where editing left off.
$ws.Range("T21").Select() | Out-Null
